$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Date" column (C) values for the export-history rows.
$ws.Range("C2").Value = (Get-Date -Year 2025 -Month 7 -Day 3 -Hour 0 -Minute 0 -Second 10)
$ws.Range("C3").Value = (Get-Date -Year 2025 -Month 8 -Day 3 -Hour 0 -Minute 0 -Second 10)
$ws.Range("C4").Value = (Get-Date -Year 2025 -Month 8 -Day 3 -Hour 0 -Minute 0 -Second 10)
$ws.Range("C5").Value = (Get-Date -Year 2025 -Month 8 -Day 3 -Hour 0 -Minute 0 -Second 10)
$ws.Range("C6").Value = (Get-Date -Year 2025 -Month 8 -Day 3 -Hour 0 -Minute 0 -Second 10)

# Update the "Food Items" text for row 4 to include the additional item.
$ws.Range("H4").Value = "Strawberry Lassi, Chicken Wrap"
